# GitHub Actions price-refresh: update the cryptos list "Price" (D) and
# "Volume(1h)" (E) columns to the newly scraped figures.
#
# D-column cells hold numeric-looking text (e.g. "213.47", or
# "27.956.02" which even has extra "thousands" dots). The source file
# stores them as plain (inline) strings, so each value we write is
# prefixed with a leading apostrophe -- Excel's standard "force text"
# marker -- to stop it from being auto-coerced into a Number cell.
# E-column values already contain padding spaces/percent signs, so they
# stay text without needing that treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.956.02'
$ws.Range("E2").Value = '  +1.62%  '
$ws.Range("D3").Value = '''1.646.83'
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''213.47'
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("D6").Value = '''0.527'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''23.40'
$ws.Range("E8").Value = '  +2.34%  '
$ws.Range("D9").Value = '''0.265'
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("D11").Value = '''0.0871'
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").Value = '''1.880.71'
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("D13").Value = '''1.639.91'
$ws.Range("E13").Value = '  +1.45%  '
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").Value = '''0.565'
$ws.Range("E15").Value = '  +2.67%  '
$ws.Range("D16").Value = '''65.57'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '''27.968.58'
$ws.Range("E17").Value = '  +1.75%  '
$ws.Range("D18").Value = '''232.61'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").Value = '''7.69'
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("D20").Value = '''0.0₃0722'
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = '''10.66'
$ws.Range("E22").Value = '  +4.51%  '
$ws.Range("D23").Value = '''4.39'
$ws.Range("E23").Value = '  +2.48%  '
$ws.Range("D24").Value = '''2.15'
$ws.Range("E24").Value = '  +3.62%  '
$ws.Range("D25").Value = '''152.29'
$ws.Range("E25").Value = '  +0.89%  '
$ws.Range("D26").Value = '''6.91'
$ws.Range("E26").Value = '  +1.05%  '
$ws.Range("D27").Value = '''15.76'
$ws.Range("E27").Value = '  +1.51%  '
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = '''1.19'
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("E32").Value = '  +2.67%  '
$ws.Range("D33").Value = '''1.443.67'
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("D34").Value = '''3.08'
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").Value = '  +1.82%  '
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").Value = '''0.888'
$ws.Range("E37").Value = '  +3.29%  '
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").Value = '''0.918'
$ws.Range("E40").Value = '  -3.03%  '
$ws.Range("D41").Value = '''69.38'
$ws.Range("E41").Value = '  +2.39%  '
$ws.Range("E42").Value = '  +3.49%  '
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("E46").Value = '  +2.27%  '
$ws.Range("E47").Value = '  +4.76%  '
$ws.Range("D48").Value = '''1.788.93'
$ws.Range("E48").Value = '  +1.61%  '
$ws.Range("D49").Value = '''88.97'
$ws.Range("E50").Value = '  -0.44%  '
